$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("material_and_methods")

# Row 4 holds the "single set of methods" values for study Eagle_et_al_2021.
# Fill in previously-blank placeholder fields ahead of manager's review.
# Order matters for how new shared-string entries get appended.
$ws.Range("AF4").Value = "Plum software used"  # age_depth_model_notes
$ws.Range("AA4").Value = "CRS"                 # excess_pb210_model (was "Plum")
$ws.Range("Z4").Value  = "mass accumulation"   # excess_pb210_rate
$ws.Range("AB4").Value = "selected intervals"  # ra226_assumption (was "each sample")
